$d = $word.ActiveDocument

# 1. Update the version number on the title page: "versión 1.0" -> "versión 1.2"
$d.Content.Find.Execute("versión 1.0", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "versión 1.2", 2)

# 2. Update the signing date: "Lima, 01 de Mayo del 2024" -> "Lima, 03 de julio del 2024"
$d.Content.Find.Execute("Lima, 01 de Mayo del 2024", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Lima, 03 de julio del 2024", 2)

# 3. Append a new revision-history row to the first table (Fecha / Versión / Descripción / Autor)
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "06/06/2024"
$newRow.Cells.Item(2).Range.Text = "1.2"
$newRow.Cells.Item(3).Range.Text = "verificación, actualización y finalización del Documento"
$newRow.Cells.Item(4).Range.Text = "Pablo Mendoza"
